# File chooser aggiunti altri impliciti
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the "implicit" estimated/actual time entries for rows that did not have them yet.
$ws.Range("J19").Value = "1gg"
$ws.Range("K19").Value = "2gg"

$ws.Range("J24").Value = "2gg"
$ws.Range("K24").Value = "3gg"

# Bump the time estimate for the Sudoku loading-interface file chooser row.
$ws.Range("K27").Value = "5gg"

# Restore the view/selection state captured in the workbook.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("K21").Select()
